$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row above row 2 (shifts existing rows down)
$ws.Rows.Item(2).Insert()

# The inserted row inherits formatting from a neighboring row; reset it back to
# the plain/default style so it matches the unstyled data rows in this sheet.
$ws.Rows.Item(2).ClearFormats()

# Fill in the new row 2 data: Japanese J League 3 match
$ws.Range("A2").Value = "Japanese J League 3"
# Force the date-looking text to stay a plain string (avoid Excel auto-converting
# it to a real date), then clear the cell format so no stray style lingers.
$ws.Range("B2").Value = "'2025-11-10"
$ws.Range("B2").ClearFormats()
$ws.Range("C2").Value = "03:00:00"
$ws.Range("D2").Value = "FC Osaka"
$ws.Range("E2").Value = "Nagano Parceiro"
$ws.Range("F2").Value = 1.02
$ws.Range("G2").Value = 1000
$ws.Range("H2").Value = 1.02
$ws.Range("I2").Value = 1000
$ws.Range("J2").Value = 1.02
$ws.Range("K2").Value = 1000
$ws.Range("L2").Value = 1.01
$ws.Range("M2").Value = 1.01
$ws.Range("N2").Value = 1.25
$ws.Range("O2").Value = 1.01
$ws.Range("P2").Value = 1.24
$ws.Range("Q2").Value = 1.02
$ws.Range("R2").Value = 1.18
$ws.Range("S2").Value = 1.02
$ws.Range("T2").Value = 1.01
$ws.Range("U2").Value = 1.01
$ws.Range("V2").Value = 1.01
$ws.Range("W2").Value = 1.01
$ws.Range("X2").Value = 1000
$ws.Range("Y2").Value = 1000
$ws.Range("Z2").Value = 1000
$ws.Range("AA2").Value = 1000
$ws.Range("AB2").Value = 1000
$ws.Range("AC2").Value = 1000
$ws.Range("AD2").Value = 1000
$ws.Range("AE2").Value = 1000
$ws.Range("AF2").Value = 1000
$ws.Range("AG2").Value = 1000
$ws.Range("AH2").Value = 1000
$ws.Range("AI2").Value = 1000
$ws.Range("AJ2").Value = 1000
$ws.Range("AK2").Value = 1000
$ws.Range("AL2").Value = 1000
$ws.Range("AM2").Value = 1000
$ws.Range("AN2").Value = 1000
$ws.Range("AO2").Value = 1000

# Update odds for the now-shifted rows (original rows 2-8 -> now rows 3-9)

# Row 3: Danish 1st Division - Esbjerg vs Middelfart
$ws.Range("F3").Value = 1.62
$ws.Range("G3").Value = 1.96
$ws.Range("H3").Value = 4
$ws.Range("I3").Value = 980
$ws.Range("J3").Value = 2.04
$ws.Range("K3").Value = 7.8
$ws.Range("Q3").Value = 1.64

# Row 4: English Premier League 2 - Div 1 - Blackburn U21 vs Derby U21
$ws.Range("K4").Value = 980

# Row 6: Argentinian Primera Division - Gimnasia La Plata vs Velez Sarsfield
$ws.Range("H6").Value = 2.46
$ws.Range("I6").Value = 2.66

# Row 8: Brazilian Serie B - Botafogo SP vs Amazonas FC
$ws.Range("F8").Value = 2.26
$ws.Range("G8").Value = 2.4
$ws.Range("I8").Value = 4.2

# Row 9: Brazilian Serie B - Chapecoense vs America MG
$ws.Range("F9").Value = 2.04
$ws.Range("J9").Value = 3.25
$ws.Range("P9").Value = 1.68
$ws.Range("Q9").Value = 2
